$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2: "ddd" -> "gg"
$ws.Range("A2").Value = "gg"

# F2 / K2 hold text digits ("0"/"1"), not numbers. A plain .Value = "1"
# would be auto-coerced to a number by Excel, so force a Text number
# format first (keeps the stored cell type as Text/string), then reset
# the cell's style back to Normal so no stray formatting is left behind.
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "1"
$ws.Range("F2").Style = "Normal"

$ws.Range("K2").NumberFormat = "@"
$ws.Range("K2").Value = "0"
$ws.Range("K2").Style = "Normal"

# F3 / K3 are real numeric cells.
$ws.Range("F3").Value = 1
$ws.Range("K3").Value = 0
